# Reorders the data rows (2-6) of Sheet1 so the "Shooting fuels debate over
# safety of Prozac for teens" record (originally row 4) moves up to row 2,
# pushing "Red Lake Reservation Readies Burial Rituals" (row 2) and
# "Everyday Hero: Jeff May" (row 3) down to rows 3 and 4, and swapping
# "Tribe Buries 3 on a Long Road to Healing" (row 6) with
# "Juggalos take issue with label as a gang | Modesto Bee" (row 5).
# Row 7 ("Rodgers Law Office") is left untouched.
#
# This is the "added one json for time bucket analysis" commit: a new
# record's worth of data lands in the right chronological slot, shuffling
# the previously-adjacent rows down/around it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Snapshot the current values (rows 2-6, columns A-E) and the
#     hyperlink target each row's column-E cell carries. ---
$rows = @(2,3,4,5,6)
$data = @{}
foreach ($r in $rows) {
    $data[$r] = @{
        A = $ws.Cells.Item($r,1).Value2
        B = $ws.Cells.Item($r,2).Value2
        C = $ws.Cells.Item($r,3).Value2
        D = $ws.Cells.Item($r,4).Value2
        E = $ws.Cells.Item($r,5).Value2
    }
}
$linkByRow = @{}
foreach ($h in $ws.Hyperlinks) {
    $linkByRow[$h.Range.Row] = $h.Address
}

# New row order: row N's new content comes from this old row number.
$order = @{
    2 = 4
    3 = 2
    4 = 3
    5 = 6
    6 = 5
}

# --- Write the cell values in their new positions. ---
foreach ($newRow in 2..6) {
    $oldRow = $order[$newRow]
    $rec = $data[$oldRow]

    $ws.Cells.Item($newRow,1).Value = $rec.A
    $ws.Cells.Item($newRow,2).Value = $rec.B
    $ws.Cells.Item($newRow,3).Value = $rec.C
    $ws.Cells.Item($newRow,4).Value = $rec.D
    $ws.Cells.Item($newRow,5).Value = $rec.E
}

# --- Re-point each row's existing hyperlink object (in place, so the
#     relationship id / cell style stay exactly as they were) to the
#     target that belongs with its row's new content. ---
foreach ($h in $ws.Hyperlinks) {
    $newRow = $h.Range.Row
    if ($order.ContainsKey($newRow)) {
        $oldRow = $order[$newRow]
        $h.Address = $linkByRow[$oldRow]
    }
}
